$wb = $excel.ActiveWorkbook

# Delete the "Desarquivamentos Pendentes" sheet
$wsDel = $wb.Worksheets.Item("Desarquivamentos Pendentes")
[void]$wsDel.Delete()

# Rename "Paineis DARQ" -> "PAINEIS DARQ"
$wsPaineis = $wb.Worksheets.Item("Paineis DARQ")
$wsPaineis.Name = "PAINEIS DARQ"

# Rename "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO"
$wsRecolhimento = $wb.Worksheets.Item("Recolhimento x Eliminacao")
$wsRecolhimento.Name = "RECOLHIMENTO X ELIMINAÇÃO"

# Restore the originally active/selected sheet (deleting a sheet can shift
# the active tab); keep "PAINEIS DARQ" selected as it was before the edit.
$wsPaineis.Activate()
